# New crime data collected - weekly CompStat update for 42nd Precinct
# Updates: report header (volume number, week-covering dates) and the
# Week-to-Date / 28-Day / Year-to-Date / Historical crime-count table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header rich text: "Volume 30   Number  35" -> "...Number  36"
# ---------------------------------------------------------------------
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "36"

# "Report Covering the Week  8/28/2023  Through  9/3/2023"
#   -> "Report Covering the Week  9/4/2023  Through  9/10/2023"
# (replace the later substring first so the earlier offset stays valid)
$c9 = $ws.Range("C9")
$c9.Characters(47, 8).Text = "9/10/2023"
$c9.Characters(27, 9).Text = "9/4/2023"

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
$ws.Range("L14").Value = -43.75
$ws.Range("N14").Value = -57.142857142857

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("L15").Value = -8

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -44.444444444444
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = -5.555555555555
$ws.Range("I16").Value = 292
$ws.Range("J16").Value = 306
$ws.Range("K16").Value = -4.575163398692
$ws.Range("L16").Value = 31.531531531531
$ws.Range("M16").Value = 56.149732620320
$ws.Range("N16").Value = -65.647058823529

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -21.052631578947
$ws.Range("F17").Value = 69
$ws.Range("G17").Value = 52
$ws.Range("H17").Value = 32.692307692307
$ws.Range("I17").Value = 497
$ws.Range("J17").Value = 466
$ws.Range("K17").Value = 6.652360515021
$ws.Range("L17").Value = 14.780600461893
$ws.Range("M17").Value = 130.092592592593
$ws.Range("N17").Value = -28.282828282828

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -29.166666666666
$ws.Range("I18").Value = 136
$ws.Range("J18").Value = 237
$ws.Range("K18").Value = -42.616033755274
$ws.Range("L18").Value = 4.615384615384
$ws.Range("M18").Value = 34.653465346534
$ws.Range("N18").Value = -81.241379310344

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -10.714285714285
$ws.Range("I19").Value = 346
$ws.Range("J19").Value = 355
$ws.Range("K19").Value = -2.535211267605
$ws.Range("L19").Value = 16.498316498316
$ws.Range("M19").Value = 108.433734939759
$ws.Range("N19").Value = 38.955823293172

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 44.444444444444
$ws.Range("F20").Value = 41
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 70.833333333333
$ws.Range("I20").Value = 323
$ws.Range("J20").Value = 209
$ws.Range("K20").Value = 54.545454545454
$ws.Range("L20").Value = 130.714285714286
$ws.Range("M20").Value = 308.860759493671
$ws.Range("N20").Value = -6.916426512968

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = -15.094339622641
$ws.Range("F21").Value = 213
$ws.Range("G21").Value = 193
$ws.Range("H21").Value = 10.362694300518
$ws.Range("I21").Value = 1626
$ws.Range("J21").Value = 1607
$ws.Range("K21").Value = 1.182327317983
$ws.Range("L21").Value = 28.741092636579
$ws.Range("M21").Value = 111.71875
$ws.Range("N21").Value = -44.769021739130

# ---------------------------------------------------------------------
# Row 22 - Transit (F/G/H flip from numeric counts to the "no data"
# placeholder text used elsewhere in this table: "0" / "***.*")
# ---------------------------------------------------------------------
$ws.Range("D30").Copy($ws.Range("F22"))
$ws.Range("F22").Text = "0"

$ws.Range("D30").Copy($ws.Range("G22"))
$ws.Range("G22").Text = "0"

$ws.Range("E30").Copy($ws.Range("H22"))
$ws.Range("H22").Text = "***.*"

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = -30
$ws.Range("F23").Value = 37
$ws.Range("H23").Value = 27.586206896551
$ws.Range("I23").Value = 286
$ws.Range("J23").Value = 247
$ws.Range("K23").Value = 15.789473684210
$ws.Range("L23").Value = 104.285714285714
$ws.Range("M23").Value = 113.432835820896

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -37.323943661971
$ws.Range("I24").Value = 823
$ws.Range("J24").Value = 895
$ws.Range("K24").Value = -8.044692737430
$ws.Range("L24").Value = 27.399380804953
$ws.Range("M24").Value = 55.576559546313

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 73.333333333333
$ws.Range("F25").Value = 80
$ws.Range("H25").Value = 2.564102564102
$ws.Range("I25").Value = 773
$ws.Range("J25").Value = 702
$ws.Range("K25").Value = 10.113960113960
$ws.Range("L25").Value = 27.980132450331
$ws.Range("M25").Value = 26.721311475409

# ---------------------------------------------------------------------
# Row 26 - UCR Rape* (C flips from placeholder text "0" to a real count)
# ---------------------------------------------------------------------
$ws.Range("F30").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1

$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 39
$ws.Range("K26").Value = -13.333333333333
$ws.Range("L26").Value = -2.5

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 69
$ws.Range("K27").Value = 30.188679245283
$ws.Range("L27").Value = 53.333333333333

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic. (C flips placeholder -> count, D/E flip
# count -> placeholder text)
# ---------------------------------------------------------------------
$ws.Range("F30").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2

$ws.Range("D30").Copy($ws.Range("D28"))
$ws.Range("D28").Text = "0"

$ws.Range("E30").Copy($ws.Range("E28"))
$ws.Range("E28").Text = "***.*"

$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 31
$ws.Range("K28").Value = 10.714285714285
$ws.Range("L28").Value = -47.457627118644
$ws.Range("M28").Value = -3.125
$ws.Range("N28").Value = -59.210526315789

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc. (C flips placeholder -> count, D/E flip
# count -> placeholder text)
# ---------------------------------------------------------------------
$ws.Range("F30").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 2

$ws.Range("D30").Copy($ws.Range("D29"))
$ws.Range("D29").Text = "0"

$ws.Range("E30").Copy($ws.Range("E29"))
$ws.Range("E29").Text = "***.*"

$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 26
$ws.Range("K29").Value = 4
$ws.Range("L29").Value = -46.938775510204
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -65.789473684210
